$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Both header stories (primary + first-page/even) carry the BTec logo
# (wp:docPr / pic:cNvPr name "image1.jpg") -> rename to "image2.jpg".
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers($i)
    if ($hdr.Exists) {
        $ishapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $ishapes.Count; $j++) {
            $shp = $ishapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}

# Both footer stories carry the Pearson Edexcel logo
# (wp:docPr / pic:cNvPr name "image2.png") -> rename to "image1.png".
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers($i)
    if ($ftr.Exists) {
        $ishapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $ishapes.Count; $j++) {
            $shp = $ishapes.Item($j)
            if ($shp.AlternativeText -like "*PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}

Write-Output "Renamed header/footer logo images."
